$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "abc"
$ws.Range("C2").Value = "fgh"
$ws.Range("C2").Select()
